$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cells that are unambiguously text (Excel will not coerce them to numbers)
$ws.Cells.Item(2, 4).Formula = '67.676.83'
$ws.Cells.Item(2, 5).Formula = '  -0.67%  '
$ws.Cells.Item(3, 4).Formula = '3.791.20'
$ws.Cells.Item(3, 5).Formula = '  +0.57%  '
$ws.Cells.Item(4, 5).Formula = '  +0.02%  '
$ws.Cells.Item(5, 5).Formula = '  +0.43%  '
$ws.Cells.Item(6, 5).Formula = '  -0.30%  '
$ws.Cells.Item(7, 4).Formula = '3.789.50'
$ws.Cells.Item(7, 5).Formula = '  +0.60%  '
$ws.Cells.Item(8, 5).Formula = '  +0.04%  '
$ws.Cells.Item(9, 5).Formula = '  +0.20%  '
$ws.Cells.Item(10, 5).Formula = '  -0.06%  '
$ws.Cells.Item(11, 5).Formula = '  -1.21%  '
$ws.Cells.Item(12, 5).Formula = '  -0.19%  '
$ws.Cells.Item(13, 5).Formula = '  -2.15%  '
$ws.Cells.Item(14, 5).Formula = '  -0.08%  '
$ws.Cells.Item(15, 4).Formula = '4.428.27'
$ws.Cells.Item(15, 5).Formula = '  +0.57%  '
$ws.Cells.Item(16, 4).Formula = '3.802.70'
$ws.Cells.Item(16, 5).Formula = '  +1.28%  '
$ws.Cells.Item(17, 5).Formula = '  +3.58%  '
$ws.Cells.Item(18, 4).Formula = '67.665.29'
$ws.Cells.Item(18, 5).Formula = '  -0.70%  '
$ws.Cells.Item(19, 5).Formula = '  +0.80%  '
$ws.Cells.Item(20, 5).Formula = '  +0.08%  '
$ws.Cells.Item(21, 5).Formula = '  -7.10%  '
$ws.Cells.Item(22, 5).Formula = '  -1.02%  '
$ws.Cells.Item(23, 5).Formula = '  +0.14%  '
$ws.Cells.Item(24, 5).Formula = '  +3.86%  '
$ws.Cells.Item(25, 5).Formula = '  -0.43%  '
$ws.Cells.Item(26, 5).Formula = '  +2.76%  '
$ws.Cells.Item(27, 5).Formula = '  -3.28%  '
$ws.Cells.Item(28, 5).Formula = '  -0.08%  '
$ws.Cells.Item(29, 5).Formula = '  -0.98%  '
$ws.Cells.Item(30, 5).Formula = '  -0.23%  '
$ws.Cells.Item(31, 5).Formula = '  +4.51%  '
$ws.Cells.Item(32, 5).Formula = '  -1.16%  '
$ws.Cells.Item(33, 5).Formula = '  -0.82%  '
$ws.Cells.Item(34, 5).Formula = '  -0.02%  '
$ws.Cells.Item(35, 5).Formula = '  -0.68%  '
$ws.Cells.Item(36, 2).Formula = 'Hedera'
$ws.Cells.Item(36, 3).Formula = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(36, 5).Formula = '  -0.36%  '
$ws.Cells.Item(37, 2).Formula = 'dogwifhat'
$ws.Cells.Item(37, 3).Formula = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(37, 5).Formula = '  -2.63%  '
$ws.Cells.Item(38, 2).Formula = 'Kaspa'
$ws.Cells.Item(38, 3).Formula = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(38, 5).Formula = '  -0.48%  '
$ws.Cells.Item(39, 2).Formula = 'Mantle'
$ws.Cells.Item(39, 3).Formula = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(39, 5).Formula = '  -0.77%  '
$ws.Cells.Item(40, 2).Formula = 'Filecoin'
$ws.Cells.Item(40, 3).Formula = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(40, 5).Formula = '  +0.06%  '
$ws.Cells.Item(41, 2).Formula = 'FirstDigitalUSD'
$ws.Cells.Item(41, 3).Formula = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(41, 5).Formula = '  -0.03%  '
$ws.Cells.Item(42, 2).Formula = 'USDe'
$ws.Cells.Item(42, 3).Formula = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(42, 5).Formula = '  +0.00%  '
$ws.Cells.Item(43, 2).Formula = 'OKB'
$ws.Cells.Item(43, 3).Formula = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(43, 5).Formula = '  +2.41%  '
$ws.Cells.Item(44, 2).Formula = 'Arweave'
$ws.Cells.Item(44, 3).Formula = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(44, 5).Formula = '  -1.18%  '
$ws.Cells.Item(45, 2).Formula = 'TheGraph'
$ws.Cells.Item(45, 3).Formula = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(45, 5).Formula = '  -0.56%  '
$ws.Cells.Item(46, 2).Formula = 'Monero'
$ws.Cells.Item(46, 3).Formula = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(46, 5).Formula = '  +3.04%  '
$ws.Cells.Item(47, 2).Formula = 'Cosmos'
$ws.Cells.Item(47, 3).Formula = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(47, 5).Formula = '  -1.61%  '
$ws.Cells.Item(48, 2).Formula = 'EnergySwap'
$ws.Cells.Item(48, 3).Formula = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 5).Formula = '  +6.90%  '
$ws.Cells.Item(49, 2).Formula = 'Bittensor'
$ws.Cells.Item(49, 3).Formula = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(49, 5).Formula = '  +0.38%  '
$ws.Cells.Item(50, 2).Formula = 'Stacks'
$ws.Cells.Item(50, 3).Formula = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(50, 5).Formula = '  -4.53%  '
$ws.Cells.Item(51, 2).Formula = 'FLOKI'
$ws.Cells.Item(51, 3).Formula = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Cells.Item(51, 5).Formula = '  +2.37%  '

# Update cells whose new text looks like a plain number, which Excel would otherwise
# auto-convert to a numeric value. Force Text format first so the literal text is kept.
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '595.72'
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '166.94'
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.521'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.159'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '6.34'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '36.04'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '18.53'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.04'
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.01'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '459.27'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '83.47'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '12.13'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.98'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.24'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '29.62'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '9.05'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.1000'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '3.36'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.137'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.993'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '5.77'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.00'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.00'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '48.04'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '43.81'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.298'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '149.76'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '8.27'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '26.95'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '389.82'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.82'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.000259'

# Normalize style back to the default (unstyled) look used by all other data cells,
# by pasting formats from a plain, unstyled text cell onto each forced cell individually.
$ws.Cells.Item(3, 2).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4122)
$ws.Cells.Item(6, 4).PasteSpecial(-4122)
$ws.Cells.Item(9, 4).PasteSpecial(-4122)
$ws.Cells.Item(10, 4).PasteSpecial(-4122)
$ws.Cells.Item(11, 4).PasteSpecial(-4122)
$ws.Cells.Item(14, 4).PasteSpecial(-4122)
$ws.Cells.Item(17, 4).PasteSpecial(-4122)
$ws.Cells.Item(19, 4).PasteSpecial(-4122)
$ws.Cells.Item(21, 4).PasteSpecial(-4122)
$ws.Cells.Item(22, 4).PasteSpecial(-4122)
$ws.Cells.Item(25, 4).PasteSpecial(-4122)
$ws.Cells.Item(26, 4).PasteSpecial(-4122)
$ws.Cells.Item(29, 4).PasteSpecial(-4122)
$ws.Cells.Item(31, 4).PasteSpecial(-4122)
$ws.Cells.Item(33, 4).PasteSpecial(-4122)
$ws.Cells.Item(35, 4).PasteSpecial(-4122)
$ws.Cells.Item(36, 4).PasteSpecial(-4122)
$ws.Cells.Item(37, 4).PasteSpecial(-4122)
$ws.Cells.Item(38, 4).PasteSpecial(-4122)
$ws.Cells.Item(39, 4).PasteSpecial(-4122)
$ws.Cells.Item(40, 4).PasteSpecial(-4122)
$ws.Cells.Item(41, 4).PasteSpecial(-4122)
$ws.Cells.Item(42, 4).PasteSpecial(-4122)
$ws.Cells.Item(43, 4).PasteSpecial(-4122)
$ws.Cells.Item(44, 4).PasteSpecial(-4122)
$ws.Cells.Item(45, 4).PasteSpecial(-4122)
$ws.Cells.Item(46, 4).PasteSpecial(-4122)
$ws.Cells.Item(47, 4).PasteSpecial(-4122)
$ws.Cells.Item(48, 4).PasteSpecial(-4122)
$ws.Cells.Item(49, 4).PasteSpecial(-4122)
$ws.Cells.Item(50, 4).PasteSpecial(-4122)
$ws.Cells.Item(51, 4).PasteSpecial(-4122)
$excel.CutCopyMode = 0
